# Remove footnote markers (e.g. " [1]", " [2]", ... " [5]") from vaccine
# description strings, and collapse embedded line breaks (Chr(10)) into a
# single space, joining wrapped labels back into one line. Also de-dupe the
# "Afluria`nQuadrivalent" shared string against the existing plain
# "Afluria Quadrivalent" string used elsewhere in the Adult Influenza sheet.

$wb = $excel.ActiveWorkbook

$map = @{}
$map["DTaP [1]"] = "DTaP "
$map["DTaP-IPV [2]"] = "DTaP-IPV "
$map["DTaP-Hep B-IPV [4]"] = "DTaP-Hep B-IPV "
$map["DTaP-IP-HI [4]"] = "DTaP-IP-HI "
$map["e-IPV [5]"] = "e-IPV "
$map["Hepatitis A Pediatric [5]"] = "Hepatitis A Pediatric "
$map["Hepatitis A-Hepatitis B 18 only [3]"] = "Hepatitis A-Hepatitis B 18 only "
$map["Hepatitis B [5]`nPediatric/Adolescent"] = "Hepatitis B  Pediatric/Adolescent"
$map["Recombivax`nHB"] = "Recombivax HB"
$map["Hib [5]"] = "Hib "
$map["HPV - Human Papillomavirus 9-valent [5]"] = "HPV - Human Papillomavirus 9-valent "
$map["MENB - Meningococcal Group B [5]"] = "MENB - Meningococcal Group B "
$map["Meningococcal Conjugate (Groups A, C, Y and W-135) [5]"] = "Meningococcal Conjugate (Groups A, C, Y and W-135) "
$map["Measles, Mumps and Rubella (MMR) [1]"] = "Measles, Mumps and Rubella (MMR) "
$map["MMR/Varicella [2]"] = "MMR/Varicella "
$map["Pneumococcal`n13-valent [5] (Pediatric)"] = "Pneumococcal 13-valent  (Pediatric)"
$map["Rotavirus, Live, Oral, Pentavalent [5]"] = "Rotavirus, Live, Oral, Pentavalent "
$map["Rotavirus, Live, Oral, Oral [5]"] = "Rotavirus, Live, Oral, Oral "
$map["Tetanus and Diphtheria Toxoids [3]"] = "Tetanus and Diphtheria Toxoids "
$map["Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis [1]"] = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$map["Varicella [5]"] = "Varicella "
$map["Hepatitis A Adult [5]"] = "Hepatitis A Adult "
$map["Hepatitis A-Hepatitis B Adult [3]"] = "Hepatitis A-Hepatitis B Adult "
$map["Hepatitis B Adult [5]"] = "Hepatitis B Adult "
$map["HPV-Human Papillomavirus 9 Valent [5]"] = "HPV-Human Papillomavirus 9 Valent "
$map["Measles, Mumps,  Rubella [1]"] = "Measles, Mumps,  Rubella "
$map["Pneumococcal`n13-valent [5]"] = "Pneumococcal 13-valent "
$map["Influenza [5]`n(Age 6 months and older)"] = "Influenza  (Age 6 months and older)"
$map["Fluzone`nQuadrivalent"] = "Fluzone Quadrivalent"
$map["Influenza [5]`n(Age 6-35 months)"] = "Influenza  (Age 6-35 months)"
$map["Fluzone`nQuadrivalent`nPediatric dose"] = "Fluzone Quadrivalent Pediatric dose"
$map["Fluarix`nQuadrivalent"] = "Fluarix Quadrivalent"
$map["FluLaval`nQuadrivalent"] = "FluLaval Quadrivalent"
$map["Influenza [5]`n(Age 4 years and older)"] = "Influenza  (Age 4 years and older)"
$map["Influenza [5]`n(Age 6 -35 months)"] = "Influenza  (Age 6 -35 months)"
$map["Influenza [5]`n(Age 36 months and older)"] = "Influenza  (Age 36 months and older)"
$map["Influenza [5]`nLive, Intranasal (Age 2-49 years)"] = "Influenza  Live, Intranasal (Age 2-49 years)"
$map["FluMist`nQuadrivalent"] = "FluMist Quadrivalent"
# De-duplicate: this wrapped label collapses to text that is already used
# elsewhere as "Afluria Quadrivalent" (no embedded line break).
$map["Afluria`nQuadrivalent"] = "Afluria Quadrivalent"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $val = $cell.Value()
            if ($val -ne $null -and $map.ContainsKey($val)) {
                $cell.Value = $map[$val]
            }
        }
    }
}
